# Update odds figures on Sheet1 to reflect the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("G13").Value = 3.5
$ws.Range("I13").Value = 1.85
$ws.Range("J13").Value = 3.75
$ws.Range("L13").Value = 2.38
$ws.Range("AC13").Value = 21
$ws.Range("AS13").Value = 101
$ws.Range("AV13").Value = 41
$ws.Range("AZ13").Value = 26
$ws.Range("BB13").Value = 67

# Row 14
$ws.Range("N14").Value = 26

# Row 17
$ws.Range("M17").Value = 1.05
$ws.Range("N17").Value = 11
$ws.Range("Q17").Value = 1.85
$ws.Range("R17").Value = 2
$ws.Range("U17").Value = 1.91
$ws.Range("V17").Value = 1.8
$ws.Range("X17").Value = 7
$ws.Range("AC17").Value = 11

# Row 22
$ws.Range("Q22").Value = 1.36

# Row 23
$ws.Range("R23").Value = 1.62

# Row 25
$ws.Range("Q25").Value = 1.95
$ws.Range("R25").Value = 1.9

# Row 29
$ws.Range("G29").Value = 3.8
$ws.Range("H29").Value = 3.5
$ws.Range("I29").Value = 1.88
$ws.Range("L29").Value = 2.47
$ws.Range("N29").Value = 7.4
$ws.Range("O29").Value = 1.29
$ws.Range("P29").Value = 3.25
$ws.Range("Q29").Value = 1.88
$ws.Range("R29").Value = 1.85
$ws.Range("S29").Value = 1.4
$ws.Range("T29").Value = 2.72
$ws.Range("W29").Value = 11
$ws.Range("X29").Value = 21
$ws.Range("AC29").Value = 7.4
$ws.Range("AF29").Value = 65
$ws.Range("AI29").Value = 9
$ws.Range("AK29").Value = 16
$ws.Range("AN29").Value = 5.6
$ws.Range("AP29").Value = 28
$ws.Range("AS29").Value = 350
$ws.Range("AT29").Value = 2.72
$ws.Range("AU29").Value = 7.3
$ws.Range("AW29").Value = 3.75
$ws.Range("AX29").Value = 9.5
$ws.Range("AY29").Value = 18.5
$ws.Range("AZ29").Value = 35

# Row 30
$ws.Range("H30").Value = 3.4
$ws.Range("I30").Value = 3.8
$ws.Range("J30").Value = 2.63
$ws.Range("L30").Value = 4
$ws.Range("M30").Value = 1.05
$ws.Range("O30").Value = 1.29
$ws.Range("Q30").Value = 1.95
$ws.Range("R30").Value = 1.9
$ws.Range("U30").Value = 1.73
$ws.Range("V30").Value = 2
$ws.Range("W30").Value = 8
$ws.Range("X30").Value = 9.5
$ws.Range("AE30").Value = 13
$ws.Range("AF30").Value = 41
$ws.Range("AG30").Value = 201
$ws.Range("AI30").Value = 19
$ws.Range("AL30").Value = 29
$ws.Range("AM30").Value = 34
$ws.Range("AY30").Value = 26

# Row 31
$ws.Range("M31").Value = 1.03
$ws.Range("O31").Value = 1.14

# Row 32
$ws.Range("M32").Value = 1.08
$ws.Range("O32").Value = 1.36

# Row 33
$ws.Range("M33").Value = 1.03
$ws.Range("O33").Value = 1.22
